# Adds two new bullet items ("Critères d'évaluation" list, numId=3) right
# after the "On doit pouvoir déclarer des nouvelles variables ..." item and
# before "Le langage doit supporter la déclaration et l'appel des fonctions...".

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find (robust to any offset drift).
$anchor = $d.Content
$found = $anchor.Find.Execute("On doit pouvoir déclarer des nouvelles variables et assigner des variables existantes.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchorParagraph = $anchor.Paragraphs(1)
$insertionPoint = $d.Range($anchorParagraph.Range.End, $anchorParagraph.Range.End)

$text1 = "Le langage doit supporter les if statements avec des comparaisons (==, ≥, <, !=, etc.) et des opérateurs logiques (&&, || et !)."
$text2 = "Le langage doit supporter les loops for et/ou les loops while, incluant les mot-clé break et continue."

# Each trailing carriage return closes off a new list paragraph that
# inherits the numPr (numId=3) / pStyle of the paragraph it follows,
# matching the bullet list the anchor paragraph belongs to.
$insertionPoint.InsertAfter($text1 + [char]13 + $text2 + [char]13)
